$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B3").Value = "Iphone"
